$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("C1").Value = "ID"
$ws.Range("D1").Value = "DESCRIPTION"
$ws.Range("E1").Value = "QTY"
$ws.Range("F1").Value = "PRICE"
$ws.Range("G1").Value = "TOTAL"

# Line items
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "4Tech keyboard black "
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 600
$ws.Range("G2").Value = 600

$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "A4Tech HS-800 headphone "
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 900
$ws.Range("G3").Value = 900

$ws.Range("C4").Value = 3
$ws.Range("D4").Value = "Asus Memo Pad Tablet "
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 7800
$ws.Range("G4").Value = 7800

$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "HP Desktop C2500 Keyboard+Mouse "
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1500
$ws.Range("G5").Value = 1500

$ws.Range("C6").Value = 5
$ws.Range("D6").Value = "Logitech B170 Wireless Mouse (Black) "
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 600
$ws.Range("G6").Value = 1200

$ws.Range("C7").Value = 6
$ws.Range("D7").Value = "Benq G2020HDA Screen. "
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1500
$ws.Range("G7").Value = 3000

$ws.Range("C8").Value = 7
$ws.Range("D8").Value = "Logitech B525 Commercial HD Webcam "
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 2000
$ws.Range("G8").Value = 2000

# Totals
$ws.Range("F9").Value = "Sub Total "
$ws.Range("G9").Value = 17000

$ws.Range("F10").Value = "GST 8% "
$ws.Range("G10").Value = 1360

$ws.Range("F11").Value = "Total "
$ws.Range("G11").Value = 18360
